# Applies the "normal_z_scores_percentages" edit:
#  1. Break/remove the (unused) external link to test.xlsx — drops
#     <externalReferences> from workbook.xml and the externalLink1 part.
#  2. Fix the last boundary label "97% to 100%" -> "97.7% to 100%"
#     (shared string used by B52:B62).
#  3. Re-enter the C-column z-score formula in two passes (C2:C33 and
#     C34:C62) so Excel records them as shared formulas, matching the
#     saved file's <f t="shared" .../> grouping.
#  4. Move the active selection from B51 to B52.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Remove the external reference to test.xlsx entirely.
$wb.BreakLink("test.xlsx", 1)

# 2. Correct the boundary label text for the last bucket. Leading "'"
#    keeps it entered as text (preserves the quote-prefix cell style the
#    original cells already carried) instead of stripping it.
$ws.Range("B52:B62").Value = "'97.7% to 100%"

# 3. Re-write the z-score formulas so they are grouped as two shared
#    formulas, split at row 34 (mirrors the authored edit).
$ws.Range("C2:C33").Formula = "=(A2-AVERAGE(A:A))/10"
$ws.Range("C34:C62").Formula = "=(A34-AVERAGE(A:A))/10"

# 4. Move the selection cursor down one row.
$ws.Range("B52").Select()
